# Auto-generated edit script: updates market-price derived columns (H-N)
# on several Leve rows across multiple sheets, per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3497.2703
$ws.Range("I64").Value = 3340
$ws.Range("J64").Value = 3604.5
$ws.Range("K64").Value = 3340
$ws.Range("L64").Value = 3604.5
$ws.Range("M64").Value = -3092
$ws.Range("N64").Value = -4100.5

$ws.Range("H67").Value = 3497.2703
$ws.Range("I67").Value = 3340
$ws.Range("J67").Value = 3604.5
$ws.Range("K67").Value = 3340
$ws.Range("L67").Value = 3604.5
$ws.Range("M67").Value = -2482
$ws.Range("N67").Value = -5320.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 221
$ws.Range("I5").Value = 293.16666
$ws.Range("J5").Value = 76.666664
$ws.Range("K5").Value = 293.16666
$ws.Range("L5").Value = 76.666664
$ws.Range("M5").Value = -181.16666
$ws.Range("N5").Value = -300.666664

$ws.Range("H16").Value = 1503701.5
$ws.Range("I16").Value = 2003268.6
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 2003268.6
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -2002981.6
$ws.Range("N16").Value = -5574

$ws.Range("H32").Value = 11633369
$ws.Range("I32").Value = 13336534
$ws.Range("K32").Value = 13336534
$ws.Range("M32").Value = -13336247

$ws.Range("H74").Value = 4208.343
$ws.Range("I74").Value = 4525.448
$ws.Range("K74").Value = 4525.448
$ws.Range("M74").Value = -3651.448

$ws.Range("H77").Value = 4208.343
$ws.Range("I77").Value = 4525.448
$ws.Range("K77").Value = 22627.24
$ws.Range("M77").Value = -18259.24

$ws.Range("H80").Value = 22125
$ws.Range("I80").Value = 29000
$ws.Range("J80").Value = 19833.334
$ws.Range("K80").Value = 29000
$ws.Range("L80").Value = 19833.334
$ws.Range("M80").Value = -28002
$ws.Range("N80").Value = -21829.334

$ws.Range("H83").Value = 22125
$ws.Range("I83").Value = 29000
$ws.Range("J83").Value = 19833.334
$ws.Range("K83").Value = 87000
$ws.Range("L83").Value = 59500.00199999999
$ws.Range("M83").Value = -82008
$ws.Range("N83").Value = -69484.00199999999

$ws.Range("H110").Value = 2395.4583
$ws.Range("I110").Value = 1998.3
$ws.Range("J110").Value = 2679.1428
$ws.Range("K110").Value = 1998.3
$ws.Range("L110").Value = 2679.1428
$ws.Range("M110").Value = 46.70000000000005
$ws.Range("N110").Value = -6769.1428

$ws.Range("H124").Value = 11583.223
$ws.Range("J124").Value = 11583.223
$ws.Range("L124").Value = 11583.223
$ws.Range("N124").Value = -21403.223

$ws.Range("H125").Value = 21630.715
$ws.Range("J125").Value = 21630.715
$ws.Range("L125").Value = 21630.715
$ws.Range("N125").Value = -31470.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 221
$ws.Range("I4").Value = 293.16666
$ws.Range("J4").Value = 76.666664
$ws.Range("K4").Value = 293.16666
$ws.Range("L4").Value = 76.666664
$ws.Range("M4").Value = -178.16666
$ws.Range("N4").Value = -306.666664

$ws.Range("H22").Value = 168.83333
$ws.Range("I22").Value = 168.83333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 168.83333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 4.166670000000011
$ws.Range("N22").ClearContents()

$ws.Range("H47").Value = 43000
$ws.Range("J47").Value = 43000
$ws.Range("L47").Value = 43000
$ws.Range("N47").Value = -44040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.666668
$ws.Range("I7").Value = 22
$ws.Range("J7").Value = 83.59999999999999
$ws.Range("K7").Value = 22
$ws.Range("L7").Value = 83.59999999999999
$ws.Range("M7").Value = 91
$ws.Range("N7").Value = -309.6

$ws.Range("H22").Value = 177.5
$ws.Range("I22").Value = 174.42857
$ws.Range("J22").Value = 184.66667
$ws.Range("K22").Value = 174.42857
$ws.Range("L22").Value = 184.66667
$ws.Range("M22").Value = 175.57143
$ws.Range("N22").Value = -884.6666700000001

$ws.Range("H62").Value = 2950.0908
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2939
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2939
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4187

$ws.Range("H65").Value = 2950.0908
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2939
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 14695
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -20935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 10001740
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 10001740
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 30005220
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -30007840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100.30769
$ws.Range("I2").Value = 102.4
$ws.Range("J2").Value = 93.333336
$ws.Range("K2").Value = 102.4
$ws.Range("L2").Value = 93.333336
$ws.Range("M2").Value = 10.59999999999999
$ws.Range("N2").Value = -319.333336

$ws.Range("H70").Value = 5066.1885
$ws.Range("I70").Value = 5097.778
$ws.Range("K70").Value = 5097.778
$ws.Range("M70").Value = -4827.778

$ws.Range("H73").Value = 5066.1885
$ws.Range("I73").Value = 5097.778
$ws.Range("K73").Value = 5097.778
$ws.Range("M73").Value = -4161.778

$ws.Range("H126").Value = 2897.5
$ws.Range("I126").Value = 1556
$ws.Range("J126").Value = 4239
$ws.Range("K126").Value = 4668
$ws.Range("L126").Value = 12717
$ws.Range("M126").Value = -2198
$ws.Range("N126").Value = -17657

$ws.Range("H132").Value = 2866.7083
$ws.Range("I132").Value = 2789.647
$ws.Range("J132").Value = 3053.8572
$ws.Range("K132").Value = 8368.940999999999
$ws.Range("L132").Value = 9161.571599999999
$ws.Range("M132").Value = -5838.940999999999
$ws.Range("N132").Value = -14221.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H127").Value = 69350
$ws.Range("J127").Value = 69350
$ws.Range("L127").Value = 69350
$ws.Range("N127").Value = -79270
